# Update quote/price figures in the "Produtos Atualizados" sheet.
# These columns (Cotação, Preço Base Reais, Preço Final) store their
# numeric-looking figures as literal TEXT (shared strings), not as
# numbers. A plain Range.Value assignment would have Excel auto-detect
# the string as a number, which changes the cell type. To keep the
# values as text (matching the original workbook's cell typing) we
# temporarily force a Text number format while writing the value, then
# clear the formatting again so the cell's style/appearance is left
# exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2 - Câmera Canon
Set-TextValue "D2" "5.64"
Set-TextValue "E2" "5635.74"
Set-TextValue "G2" "7890.04"

# Row 3 - Carro Renault
Set-TextValue "D3" "6.40"
Set-TextValue "E3" "28812.75"
Set-TextValue "G3" "57625.49"

# Row 4 - Notebook Dell
Set-TextValue "D4" "5.64"
Set-TextValue "E4" "5072.16"
Set-TextValue "G4" "8622.68"

# Row 5 - IPhone
Set-TextValue "D5" "5.64"
Set-TextValue "E5" "4503.00"
Set-TextValue "G5" "7655.11"

# Row 6 - Carro Fiat
Set-TextValue "D6" "6.40"
Set-TextValue "E6" "19208.50"
Set-TextValue "G6" "36496.14"

# Row 7 - Celular Xiaomi
Set-TextValue "D7" "5.64"
Set-TextValue "E7" "2707.89"
Set-TextValue "G7" "5415.78"

# Row 8 - Joia 20g
Set-TextValue "D8" "325.48"
Set-TextValue "E8" "6509.60"
Set-TextValue "G8" "7486.04"
